{"js": "// The published HTML->docx export appended a trailing blank paragraph plus\n// the site's \"Ver no Jupiter...\" navigation line and the Jekyll footer\n// (\"\u00a9 2020 . Contact: ...\"). This rebuild of the site dropped that footer\n// block; remove the same three paragraphs that immediately follow the last\n// bibliography entry (\"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear. S\u00e3o Paulo: Thomson,\n// 2007.\") while keeping the blank paragraph (and page-break paragraph) that\n// come after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the bibliography's last line so the deletion is anchored to the\n// surrounding content rather than a hard-coded index.\nconst anchorText = \"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear. S\u00e3o Paulo: Thomson, 2007.\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the bibliography anchor paragraph.\");\n}\n\n// The three paragraphs to drop: the blank spacer, the \"Ver no Jupiter...\"\n// line, and the \"\u00a9 2020 ...\" footer line.\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 3; i++) {\n  if (i >= items.length) break;\n  toDelete.push(items[i]);\n}\n\nif (\n  toDelete.length !== 3 ||\n  toDelete[0].text !== \"\" ||\n  toDelete[1].text !== \"Ver no Jupiter Salvar em pdf Salvar em docx\" ||\n  toDelete[2].text !==\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n) {\n  throw new Error(\"Unexpected paragraph content near the footer; aborting to avoid deleting the wrong text.\");\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# The published HTML->docx export appended a trailing blank paragraph plus\n# the site's \"Ver no Jupiter...\" navigation line and the Jekyll footer\n# (\"(c) 2020 . Contact: ...\"). This rebuild of the site dropped that footer\n# block; remove the same three paragraphs that immediately follow the last\n# bibliography entry (\"... Sao Paulo: Thomson, 2007.\") while keeping the\n# blank paragraph (and page-break paragraph) that come after it.\n\n$d = $word.ActiveDocument\n\n# Anchor on the last bibliography line using an ASCII-only substring so the\n# search is not sensitive to accented-character encoding.\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"Thomson, 2007\")\nif (-not $found) {\n    throw \"Could not find the bibliography anchor text ('Thomson, 2007').\"\n}\n\n$anchorPara = $searchRange.Paragraphs(1)\n\n# The three paragraphs to drop: the blank spacer, the \"Ver no Jupiter...\"\n# line, and the footer (\"... Contact: ... Jekyll ...\") line.\n$p1 = $anchorPara.Next()\n$p2 = $p1.Next()\n$p3 = $p2.Next()\n\n$t1 = $p1.Range.Text.Trim()\n$t2 = $p2.Range.Text.Trim()\n$t3 = $p3.Range.Text\n\nif ($t1.Length -ne 0) {\n    throw \"Unexpected content in blank spacer paragraph; aborting to avoid deleting the wrong text.\"\n}\nif ($t2 -ne \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    throw \"Unexpected content in 'Ver no Jupiter' paragraph; aborting to avoid deleting the wrong text.\"\n}\nif (-not $t3.ToLower().Contains(\"jekyll\")) {\n    throw \"Unexpected content in footer paragraph; aborting to avoid deleting the wrong text.\"\n}\n\n$deleteRange = $d.Range($p1.Range.Start, $p3.Range.End)\n$deleteRange.Delete()\n"}
